$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.305.74'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '1.931.47'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9950'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7673'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +6.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '248.44'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9983'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.45'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3222'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07095'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7902'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08001'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.20%  '
$ws.Range('D13').Value = '1.932.68'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.378'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.75'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.75'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.35%  '
$ws.Range('D17').Value = '30.315.50'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.70'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008026'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.800'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('D21').Value = '2.190.03'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9991'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9966'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.828'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.599'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.74'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1360'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.78%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.322'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.12'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.363'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.442'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.145'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05167'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.295'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7531'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.769'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01965'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.798'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.46'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.423'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4519'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.994'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9993'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8364'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.79'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.558'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.814'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '991.80'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +12.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.61'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1197'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.82%  '
